$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label text is entered in this order (Final counts, Accuracy, Total
# cells) so new shared-string entries land in the same order the original
# author typed them in.
$ws.Range("A39").Value = "Final counts"
$ws.Range("A41").Value = "Accuracy"
$ws.Range("A38").Value = "Total cells"

# Row 38: "Total cells" row
$ws.Range("B38").Value = 26
$ws.Range("C38").Value = 26
$ws.Range("D38").Value = 26
$ws.Range("E38").Value = 26
$ws.Range("F38").Value = 26

# Row 39: "Final counts" row
$ws.Range("B39").Value = 7
$ws.Range("C39").Value = 21
$ws.Range("D39").Value = 9
$ws.Range("E39").Value = 11
$ws.Range("F39").Value = 21

# Row 40: blank spacer row (still styled like the rest)
$ws.Range("B40:F40").Value = ""

# Row 41: "Accuracy" row with formulas
$ws.Range("B41").Formula = "=B39/B38*100"
$ws.Range("C41").Formula = "=C39/C38*100"
$ws.Range("D41").Formula = "=D39/D38*100"
$ws.Range("E41").Formula = "=E39/E38*100"
$ws.Range("F41").Formula = "=F39/F38*100"

# Apply the centered style (style index 1 in the original workbook) used
# throughout this table to the newly added cells.
$ws.Range("B38:F41").HorizontalAlignment = -4108

# Update the view: selection moved to K15, and the frozen/top-left cell reset
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("K15").Select()
